$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values with new TPM-derived figures
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.133054
$ws.Range("H2").Value = 0.399162
$ws.Range("Q2").Value = 0.007481715122666668
$ws.Range("R2").Value = 0.067335436104
